$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales")

# Insert a new row at position 25, shifting existing rows 25-69 down to 26-70
$ws.Rows.Item(25).Insert(-4121)

# Carry over the formatting from the row above (row 24), matching Excel's
# default "insert row" behaviour of copying the format from the row above.
$ws.Range("A24:I24").Copy()
$ws.Range("A25:I25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the newly inserted row with the new sale record
$ws.Range("A25").Value = "NEW  DELUXE CO"
$ws.Range("B25").Value = "Hesseragatta"
$ws.Range("C25").Value = 125
$ws.Range("D25").Value = 45278
$ws.Range("E25").Value = 12400
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

# Update the view: scroll to A14, select H25 (the cell clicked after inserting)
$ws.Application.GoTo($ws.Range("A1"), $true)
$ws.Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H25").Select()

# Make the Sales sheet the active/selected tab (as last active sheet in the file)
$ws.Activate()

Write-Output "done"
